# ---------------------------------------------------------------------------
# TMTI0055389_EditExistingOppEngToNewCFJobType.xlsx
# "Merge - Opp Test Data, ENg Detail, Add Counterparty - 10 Oct 2025"
#
# Summary of the edit:
#  - CAOUsers gains a Group column and new CAO / SystemAdmin rows
#    (Gemma Hardy is removed).
#  - ModuleName gains an "Engagements" row.
#  - AddOpportunity: JobType value for the sample row changes from
#    "ESOP Corporate Finance" to "Directs".
#  - A brand-new "NewEngJobTypes" sheet is inserted (Engagement old/new
#    job-type conversion table).
#  - "NewJobTypes" is renamed to "NewOppJobTypes" and its conversion table
#    is replaced with the Opportunity old->new JobType mapping.
#  - "AddContact" is unchanged in content.
#  - "Engagement" sample/reference data is rebuilt with the full list of
#    new CF job types.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) StandardUsers - cosmetic only (column width / selection), no data change
# ---------------------------------------------------------------------------
$wsStd = $wb.Worksheets.Item("StandardUsers")
$wsStd.Columns.Item(1).ColumnWidth = 13.140625
[void]$wsStd.Range("F11").Select()

# ---------------------------------------------------------------------------
# 2) CAOUsers - new Group column, new roster
# ---------------------------------------------------------------------------
$wsCao = $wb.Worksheets.Item("CAOUsers")

$wsCao.Range("A1").Value = "User"
$wsCao.Range("B1").Value = "Profile"
$wsCao.Range("C1").Value = "Group"
$wsCao.Range("A1:C1").Font.Bold = $true

$wsCao.Range("A2").Value = "Giselle Segura"
$wsCao.Range("B2").Value = "CAO"
$wsCao.Range("C2").Value = "Conversion CF CS"

$wsCao.Range("A3").Value = "Ajay Nair"
$wsCao.Range("B3").Value = "SystemAdmin"

$wsCao.Columns.Item(1).AutoFit() | Out-Null
$wsCao.Columns.Item(2).AutoFit() | Out-Null
$wsCao.Columns.Item(3).AutoFit() | Out-Null
[void]$wsCao.Range("D5").Select()

# ---------------------------------------------------------------------------
# 3) AppName - cosmetic only, no data change
# ---------------------------------------------------------------------------
$wsApp = $wb.Worksheets.Item("AppName")
[void]$wsApp.Range("H30").Select()

# ---------------------------------------------------------------------------
# 4) ModuleName - add "Engagements" row
# ---------------------------------------------------------------------------
$wsMod = $wb.Worksheets.Item("ModuleName")
$wsMod.Range("A3").Value = "Engagements"
$wsMod.Columns.Item(1).AutoFit() | Out-Null
[void]$wsMod.Range("A7").Select()

# ---------------------------------------------------------------------------
# 5) AddOpportunity - JobType sample value changes, clear stray style on D2
# ---------------------------------------------------------------------------
$wsOpp = $wb.Worksheets.Item("AddOpportunity")
$wsOpp.Range("C2").Value = "Directs"
$wsOpp.Range("D2").Style = "Normal"
$wsOpp.Range("D2").Value = "BUS - Business Services"
[void]$wsOpp.Range("D6").Select()

# ---------------------------------------------------------------------------
# 6) Insert brand-new "NewEngJobTypes" sheet before "NewJobTypes"
# ---------------------------------------------------------------------------
$wsOldOppTypes = $wb.Worksheets.Item("NewJobTypes")
$wsEngTypes = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsOldOppTypes)
$wsEngTypes.Name = "NewEngJobTypes"

$wsEngTypes.Range("A1").Value = "EngName"
$wsEngTypes.Range("B1").Value = "EngNumber"
$wsEngTypes.Range("C1").Value = "OriginalValue"
$wsEngTypes.Range("D1").Value = "NewValue"
$wsEngTypes.Range("A1:D1").Font.Bold = $true

$wsEngTypes.Range("B2:B4").NumberFormat = "@"

$wsEngTypes.Range("A2").Value = "Project Grizzly + Omaha"
$wsEngTypes.Range("B2").Value = "100022"
$wsEngTypes.Range("C2").Value = "Debt Capital Markets"
$wsEngTypes.Range("D2").Value = "Debt Financing"
$wsEngTypes.Range("E2").Value = 100022

$wsEngTypes.Range("A3").Value = "Project Ergon V"
$wsEngTypes.Range("B3").Value = "100864"
$wsEngTypes.Range("C3").Value = "Private Funds: Primary Advisory"
$wsEngTypes.Range("D3").Value = "Primary Capital Advisory"
$wsEngTypes.Range("E3").Value = 100864

$wsEngTypes.Range("A4").Value = "Project Diablo - Financing"
$wsEngTypes.Range("B4").Value = "101371"
$wsEngTypes.Range("C4").Value = "Equity Capital Markets"
$wsEngTypes.Range("D4").Value = "Equity Placements"
$wsEngTypes.Range("E4").Value = 101371

$wsEngTypes.Columns.Item(1).AutoFit() | Out-Null
$wsEngTypes.Columns.Item(2).AutoFit() | Out-Null
$wsEngTypes.Columns.Item(3).AutoFit() | Out-Null
$wsEngTypes.Columns.Item(4).AutoFit() | Out-Null
[void]$wsEngTypes.Range("C9").Select()

# ---------------------------------------------------------------------------
# 7) Rename "NewJobTypes" -> "NewOppJobTypes" and replace its contents with
#    the Opportunity old -> new JobType conversion table
# ---------------------------------------------------------------------------
$wsOppTypes = $wb.Worksheets.Item("NewJobTypes")
$wsOppTypes.Name = "NewOppJobTypes"

$wsOppTypes.Cells.Clear() | Out-Null

$wsOppTypes.Range("A1").Value = "JobType"
$wsOppTypes.Range("B1").Value = "JobType"
$wsOppTypes.Range("A1:B1").Font.Bold = $true

$oppRows = @(
    @("Buyside & Financing Advisory", "Debt Financing"),
    @("Discretionary CS Advisory", "Debt Financing"),
    @("Debt Financing", "Equity Placements"),
    @("Equity Placements", "Lender Education"),
    @("Public Underwriting", "Liability Management"),
    @("Directs", "Public Underwriting"),
    @("GP Advisory", "Financial Asset Sales"),
    @("GP Stake Sale", "Directs"),
    @("GP-Led Secondaries", "GP Advisory"),
    @("LP-Led Secondaries", ""),
    @("Primary Capital Advisory", "Primary Capital Advisory"),
    @("Financial Asset Sale", "GP-Led Secondaries OR LP-led Secondaries"),
    @("Lender Education", ""),
    @("Liability Management", "")
)

$r = 2
foreach ($row in $oppRows) {
    $wsOppTypes.Cells.Item($r, 1).Value = $row[0]
    if ($row[1] -ne "") {
        $wsOppTypes.Cells.Item($r, 2).Value = $row[1]
    }
    $r++
}

$wsOppTypes.Columns.Item(1).AutoFit() | Out-Null
$wsOppTypes.Columns.Item(2).AutoFit() | Out-Null
[void]$wsOppTypes.Range("A2:A15").Select()

# ---------------------------------------------------------------------------
# 8) AddContact - unchanged content, nothing to do beyond view/selection
# ---------------------------------------------------------------------------
$wsContact = $wb.Worksheets.Item("AddContact")
[void]$wsContact.Range("L30").Select()

# ---------------------------------------------------------------------------
# 9) Engagement - rebuild sample data with the full new CF job-type list
# ---------------------------------------------------------------------------
$wsEng = $wb.Worksheets.Item("Engagement")
$wsEng.Cells.Clear() | Out-Null

$wsEng.Range("A1").Value = "JobType"
$wsEng.Range("B1").Value = "RecordType"
$wsEng.Range("A1:B1").Font.Bold = $true

$engRows = @(
    @("Buyside & Financing Advisory", "Capital Solutions"),
    @("Discretionary CS Advisory", ""),
    @("Debt Financing", ""),
    @("Equity Placements", ""),
    @("Public Underwriting", ""),
    @("Directs", ""),
    @("GP Advisory", ""),
    @("GP Stake Sale", ""),
    @("GP-Led Secondaries", ""),
    @("LP-Led Secondaries", ""),
    @("Primary Capital Advisory", ""),
    @("Financial Asset Sale", ""),
    @("Lender Education", ""),
    @("Liability Management", "")
)

$r = 2
foreach ($row in $engRows) {
    $wsEng.Cells.Item($r, 1).Value = $row[0]
    if ($row[1] -ne "") {
        $wsEng.Cells.Item($r, 2).Value = $row[1]
    }
    $r++
}

$wsEng.Columns.Item(1).AutoFit() | Out-Null
$wsEng.Columns.Item(2).AutoFit() | Out-Null

# Engagement is the sheet that ends up active / tab-selected in the source
$wsEng.Activate()
[void]$wsEng.Range("J28").Select()

# ---------------------------------------------------------------------------
# Window layout - scroll the tab strip to AppName, keep Engagement active
# ---------------------------------------------------------------------------
$wb.Windows.Item(1).ScrollWorkbookTabs(3) | Out-Null
